$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Typography": add a new "Tiny" typography row (row 8)
# ---------------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

$typo.Range("B8").Value = "Tiny"
$typo.Range("C8").Value = "trebuc.ttf"
$typo.Range("D8").Value = 16
$typo.Range("E8").Value = 4
# Columns B:E on this sheet default to column style index 1; reset these new
# cells back to the workbook's Normal style so they match their unstyled
# neighbours (rows 4-7 carry no explicit cell style either).
$typo.Range("B8:E8").Style = "Normal"

$typo.Range("F8").Value = "?"
$typo.Range("G8").Value = ',.+-*/!@#$%s^&*()_{}\`~<>?;:''" |<>'
$typo.Range("H8").Value = "0-9,a-z,A-Z"

# ---------------------------------------------------------------------------
# Sheet "Translation": rename the wildcard-size text ids and tidy their
# placeholder text, then add the new "Goto" related rows.
# ---------------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

$tr.Range("B18").Value = "Medium"
$tr.Range("F18").Value = " "

$tr.Range("B19").Value = "Small"
$tr.Range("F19").Value = " "

$tr.Range("B20").Value = "Large"
$tr.Range("F20").Value = " "

$tr.Range("B72").Value = "Tiny"
$tr.Range("C72").Value = "Tiny"
$tr.Range("D72").Value = "Left"
$tr.Range("E72").Value = "LTR"
$tr.Range("F72").Value = " "

$tr.Range("B73").Value = "SingleUseId70"
$tr.Range("C73").Value = "Small"
$tr.Range("D73").Value = "Left"
$tr.Range("E73").Value = "LTR"
$tr.Range("F73").Value = "New Text"

$tr.Range("B74").Value = "SingleUseId71"
$tr.Range("C74").Value = "Small"
$tr.Range("D74").Value = "Center"
$tr.Range("E74").Value = "LTR"
$tr.Range("F74").Value = "Go To"
